# ----------------------------------------------------------------------------
# Regenerate the per-subject stimulus list for the "living_rooms" memory block.
# The experiment now draws from a pool of 20 distinct stimulus sets (block_total
# goes from 3 -> 6) that get duplicated across the 1000-subject run, so every row
# in this sheet is re-pointed at its new stimulus/trial/rating data.
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 6  # C2 (block_total)
$ws.Cells.Item(2, 6).Value = 203  # F2 (trial_total)
$ws.Cells.Item(2, 9).Value = $null  # I2 (cond_cat)
$ws.Cells.Item(2, 10).Value = 'new'  # J2 (cond_mem)
$ws.Cells.Item(2, 11).Value = 'f'  # K2 (correct_answer)
$ws.Cells.Item(2, 12).Value = 'stimuli/img_5jp4f.png'  # L2 (stimulus)
$ws.Cells.Item(2, 13).Value = 84.85714285714286  # M2 (conceptual)
$ws.Cells.Item(2, 14).Value = 67.83333333333333  # N2 (perceptual)
$ws.Cells.Item(2, 15).Value = 76.3452380952381  # O2 (typicality)
$ws.Cells.Item(2, 16).Value = 42  # P2 (n)
$ws.Cells.Item(2, 17).Value = 9  # Q2 (p_typicality)
$ws.Cells.Item(2, 18).Value = 9  # R2 (p_conceptual)
$ws.Cells.Item(2, 19).Value = 9  # S2 (p_perceptual)
$ws.Cells.Item(2, 20).Value = 8  # T2 (r_typicality)
$ws.Cells.Item(2, 21).Value = 8  # U2 (r_conceptual)
$ws.Cells.Item(2, 22).Value = 9  # V2 (r_perceptual)

# Row 3
$ws.Cells.Item(3, 3).Value = 6  # C3 (block_total)
$ws.Cells.Item(3, 6).Value = 204  # F3 (trial_total)
$ws.Cells.Item(3, 9).Value = 'target'  # I3 (cond_cat)
$ws.Cells.Item(3, 10).Value = 'old'  # J3 (cond_mem)
$ws.Cells.Item(3, 11).Value = 'j'  # K3 (correct_answer)
$ws.Cells.Item(3, 12).Value = 'stimuli/img_o30wb.png'  # L3 (stimulus)
$ws.Cells.Item(3, 13).Value = 81.06666666666666  # M3 (conceptual)
$ws.Cells.Item(3, 14).Value = 65.37777777777778  # N3 (perceptual)
$ws.Cells.Item(3, 15).Value = 73.22222222222223  # O3 (typicality)
$ws.Cells.Item(3, 16).Value = 45  # P3 (n)
$ws.Cells.Item(3, 17).Value = 8  # Q3 (p_typicality)
$ws.Cells.Item(3, 18).Value = 8  # R3 (p_conceptual)
$ws.Cells.Item(3, 19).Value = 8  # S3 (p_perceptual)
$ws.Cells.Item(3, 20).Value = 8  # T3 (r_typicality)
$ws.Cells.Item(3, 21).Value = 8  # U3 (r_conceptual)
$ws.Cells.Item(3, 22).Value = 8  # V3 (r_perceptual)

# Row 4
$ws.Cells.Item(4, 3).Value = 6  # C4 (block_total)
$ws.Cells.Item(4, 6).Value = 205  # F4 (trial_total)
$ws.Cells.Item(4, 9).Value = $null  # I4 (cond_cat)
$ws.Cells.Item(4, 10).Value = 'new'  # J4 (cond_mem)
$ws.Cells.Item(4, 11).Value = 'f'  # K4 (correct_answer)
$ws.Cells.Item(4, 12).Value = 'stimuli/img_6a0hu.png'  # L4 (stimulus)
$ws.Cells.Item(4, 13).Value = 61.275  # M4 (conceptual)
$ws.Cells.Item(4, 14).Value = 42.025  # N4 (perceptual)
$ws.Cells.Item(4, 15).Value = 51.65  # O4 (typicality)
$ws.Cells.Item(4, 16).Value = 40  # P4 (n)
$ws.Cells.Item(4, 17).Value = 4  # Q4 (p_typicality)
$ws.Cells.Item(4, 18).Value = 4  # R4 (p_conceptual)
$ws.Cells.Item(4, 19).Value = 4  # S4 (p_perceptual)
$ws.Cells.Item(4, 20).Value = 5  # T4 (r_typicality)
$ws.Cells.Item(4, 21).Value = 4  # U4 (r_conceptual)
$ws.Cells.Item(4, 22).Value = 5  # V4 (r_perceptual)

# Row 5
$ws.Cells.Item(5, 3).Value = 6  # C5 (block_total)
$ws.Cells.Item(5, 6).Value = 206  # F5 (trial_total)
$ws.Cells.Item(5, 9).Value = 'target'  # I5 (cond_cat)
$ws.Cells.Item(5, 10).Value = 'old'  # J5 (cond_mem)
$ws.Cells.Item(5, 11).Value = 'j'  # K5 (correct_answer)
$ws.Cells.Item(5, 12).Value = 'stimuli/img_9bkl9.png'  # L5 (stimulus)
$ws.Cells.Item(5, 13).Value = 46.62162162162162  # M5 (conceptual)
$ws.Cells.Item(5, 14).Value = 34.27027027027027  # N5 (perceptual)
$ws.Cells.Item(5, 15).Value = 40.44594594594595  # O5 (typicality)
$ws.Cells.Item(5, 16).Value = 37  # P5 (n)
$ws.Cells.Item(5, 17).Value = 3  # Q5 (p_typicality)
$ws.Cells.Item(5, 18).Value = 3  # R5 (p_conceptual)
$ws.Cells.Item(5, 19).Value = 3  # S5 (p_perceptual)
$ws.Cells.Item(5, 21).Value = 3  # U5 (r_conceptual)

# Row 6
$ws.Cells.Item(6, 3).Value = 6  # C6 (block_total)
$ws.Cells.Item(6, 6).Value = 207  # F6 (trial_total)
$ws.Cells.Item(6, 12).Value = 'stimuli/img_j4ttn.png'  # L6 (stimulus)
$ws.Cells.Item(6, 13).Value = 12.61904761904762  # M6 (conceptual)
$ws.Cells.Item(6, 14).Value = 11.42857142857143  # N6 (perceptual)
$ws.Cells.Item(6, 15).Value = 12.02380952380952  # O6 (typicality)
$ws.Cells.Item(6, 17).Value = 1  # Q6 (p_typicality)
$ws.Cells.Item(6, 18).Value = 1  # R6 (p_conceptual)
$ws.Cells.Item(6, 19).Value = 1  # S6 (p_perceptual)
$ws.Cells.Item(6, 20).Value = 1  # T6 (r_typicality)
$ws.Cells.Item(6, 21).Value = 1  # U6 (r_conceptual)
$ws.Cells.Item(6, 22).Value = 1  # V6 (r_perceptual)

# Row 7
$ws.Cells.Item(7, 3).Value = 6  # C7 (block_total)
$ws.Cells.Item(7, 6).Value = 208  # F7 (trial_total)
$ws.Cells.Item(7, 9).Value = $null  # I7 (cond_cat)
$ws.Cells.Item(7, 10).Value = 'new'  # J7 (cond_mem)
$ws.Cells.Item(7, 11).Value = 'f'  # K7 (correct_answer)
$ws.Cells.Item(7, 12).Value = 'stimuli/img_16kib.png'  # L7 (stimulus)
$ws.Cells.Item(7, 13).Value = 80.97727272727273  # M7 (conceptual)
$ws.Cells.Item(7, 14).Value = 61.11363636363637  # N7 (perceptual)
$ws.Cells.Item(7, 15).Value = 71.04545454545455  # O7 (typicality)
$ws.Cells.Item(7, 16).Value = 44  # P7 (n)
$ws.Cells.Item(7, 17).Value = 8  # Q7 (p_typicality)
$ws.Cells.Item(7, 18).Value = 8  # R7 (p_conceptual)
$ws.Cells.Item(7, 19).Value = 8  # S7 (p_perceptual)
$ws.Cells.Item(7, 20).Value = 7  # T7 (r_typicality)
$ws.Cells.Item(7, 21).Value = 7  # U7 (r_conceptual)
$ws.Cells.Item(7, 22).Value = 7  # V7 (r_perceptual)

# Row 8
$ws.Cells.Item(8, 3).Value = 6  # C8 (block_total)
$ws.Cells.Item(8, 6).Value = 209  # F8 (trial_total)
$ws.Cells.Item(8, 9).Value = $null  # I8 (cond_cat)
$ws.Cells.Item(8, 10).Value = 'new'  # J8 (cond_mem)
$ws.Cells.Item(8, 11).Value = 'f'  # K8 (correct_answer)
$ws.Cells.Item(8, 12).Value = 'stimuli/img_5jy9c.png'  # L8 (stimulus)
$ws.Cells.Item(8, 13).Value = 87.37209302325581  # M8 (conceptual)
$ws.Cells.Item(8, 14).Value = 79.18604651162791  # N8 (perceptual)
$ws.Cells.Item(8, 15).Value = 83.27906976744185  # O8 (typicality)
$ws.Cells.Item(8, 16).Value = 43  # P8 (n)
$ws.Cells.Item(8, 17).Value = 10  # Q8 (p_typicality)
$ws.Cells.Item(8, 18).Value = 10  # R8 (p_conceptual)
$ws.Cells.Item(8, 19).Value = 10  # S8 (p_perceptual)
$ws.Cells.Item(8, 20).Value = 10  # T8 (r_typicality)
$ws.Cells.Item(8, 21).Value = 9  # U8 (r_conceptual)
$ws.Cells.Item(8, 22).Value = 10  # V8 (r_perceptual)

# Row 9
$ws.Cells.Item(9, 3).Value = 6  # C9 (block_total)
$ws.Cells.Item(9, 6).Value = 210  # F9 (trial_total)
$ws.Cells.Item(9, 9).Value = $null  # I9 (cond_cat)
$ws.Cells.Item(9, 10).Value = 'new'  # J9 (cond_mem)
$ws.Cells.Item(9, 11).Value = 'f'  # K9 (correct_answer)
$ws.Cells.Item(9, 12).Value = 'stimuli/img_z4jxm.png'  # L9 (stimulus)
$ws.Cells.Item(9, 13).Value = 88.30952380952381  # M9 (conceptual)
$ws.Cells.Item(9, 14).Value = 72.64285714285714  # N9 (perceptual)
$ws.Cells.Item(9, 15).Value = 80.47619047619048  # O9 (typicality)
$ws.Cells.Item(9, 16).Value = 42  # P9 (n)
$ws.Cells.Item(9, 17).Value = 10  # Q9 (p_typicality)
$ws.Cells.Item(9, 18).Value = 10  # R9 (p_conceptual)
$ws.Cells.Item(9, 19).Value = 10  # S9 (p_perceptual)
$ws.Cells.Item(9, 20).Value = 10  # T9 (r_typicality)
$ws.Cells.Item(9, 21).Value = 10  # U9 (r_conceptual)
$ws.Cells.Item(9, 22).Value = 10  # V9 (r_perceptual)

# Row 10
$ws.Cells.Item(10, 3).Value = 6  # C10 (block_total)
$ws.Cells.Item(10, 6).Value = 211  # F10 (trial_total)
$ws.Cells.Item(10, 9).Value = $null  # I10 (cond_cat)
$ws.Cells.Item(10, 10).Value = 'new'  # J10 (cond_mem)
$ws.Cells.Item(10, 11).Value = 'f'  # K10 (correct_answer)
$ws.Cells.Item(10, 12).Value = 'stimuli/img_pdzf1.png'  # L10 (stimulus)
$ws.Cells.Item(10, 13).Value = 86.23913043478261  # M10 (conceptual)
$ws.Cells.Item(10, 14).Value = 67.17391304347827  # N10 (perceptual)
$ws.Cells.Item(10, 15).Value = 76.70652173913044  # O10 (typicality)
$ws.Cells.Item(10, 16).Value = 46  # P10 (n)
$ws.Cells.Item(10, 17).Value = 9  # Q10 (p_typicality)
$ws.Cells.Item(10, 18).Value = 9  # R10 (p_conceptual)
$ws.Cells.Item(10, 19).Value = 9  # S10 (p_perceptual)
$ws.Cells.Item(10, 20).Value = 9  # T10 (r_typicality)
$ws.Cells.Item(10, 21).Value = 9  # U10 (r_conceptual)
$ws.Cells.Item(10, 22).Value = 8  # V10 (r_perceptual)

# Row 11
$ws.Cells.Item(11, 3).Value = 6  # C11 (block_total)
$ws.Cells.Item(11, 6).Value = 212  # F11 (trial_total)
$ws.Cells.Item(11, 9).Value = $null  # I11 (cond_cat)
$ws.Cells.Item(11, 10).Value = 'new'  # J11 (cond_mem)
$ws.Cells.Item(11, 11).Value = 'f'  # K11 (correct_answer)
$ws.Cells.Item(11, 12).Value = 'stimuli/img_jkm86.png'  # L11 (stimulus)
$ws.Cells.Item(11, 13).Value = 58.32558139534883  # M11 (conceptual)
$ws.Cells.Item(11, 14).Value = 38.65116279069768  # N11 (perceptual)
$ws.Cells.Item(11, 15).Value = 48.48837209302326  # O11 (typicality)
$ws.Cells.Item(11, 16).Value = 43  # P11 (n)
$ws.Cells.Item(11, 17).Value = 4  # Q11 (p_typicality)
$ws.Cells.Item(11, 18).Value = 4  # R11 (p_conceptual)
$ws.Cells.Item(11, 19).Value = 4  # S11 (p_perceptual)
$ws.Cells.Item(11, 20).Value = 4  # T11 (r_typicality)
$ws.Cells.Item(11, 21).Value = 4  # U11 (r_conceptual)
$ws.Cells.Item(11, 22).Value = 4  # V11 (r_perceptual)

# Row 12
$ws.Cells.Item(12, 3).Value = 6  # C12 (block_total)
$ws.Cells.Item(12, 6).Value = 213  # F12 (trial_total)
$ws.Cells.Item(12, 9).Value = $null  # I12 (cond_cat)
$ws.Cells.Item(12, 10).Value = 'new'  # J12 (cond_mem)
$ws.Cells.Item(12, 11).Value = 'f'  # K12 (correct_answer)
$ws.Cells.Item(12, 12).Value = 'stimuli/img_vgh2g.png'  # L12 (stimulus)
$ws.Cells.Item(12, 13).Value = 93.81395348837209  # M12 (conceptual)
$ws.Cells.Item(12, 14).Value = 78.27906976744185  # N12 (perceptual)
$ws.Cells.Item(12, 15).Value = 86.04651162790697  # O12 (typicality)
$ws.Cells.Item(12, 16).Value = 43  # P12 (n)
$ws.Cells.Item(12, 17).Value = 10  # Q12 (p_typicality)
$ws.Cells.Item(12, 18).Value = 10  # R12 (p_conceptual)
$ws.Cells.Item(12, 19).Value = 10  # S12 (p_perceptual)
$ws.Cells.Item(12, 20).Value = 10  # T12 (r_typicality)
$ws.Cells.Item(12, 21).Value = 10  # U12 (r_conceptual)
$ws.Cells.Item(12, 22).Value = 10  # V12 (r_perceptual)

# Row 13
$ws.Cells.Item(13, 3).Value = 6  # C13 (block_total)
$ws.Cells.Item(13, 6).Value = 214  # F13 (trial_total)
$ws.Cells.Item(13, 9).Value = $null  # I13 (cond_cat)
$ws.Cells.Item(13, 10).Value = 'new'  # J13 (cond_mem)
$ws.Cells.Item(13, 11).Value = 'f'  # K13 (correct_answer)
$ws.Cells.Item(13, 12).Value = 'stimuli/img_i6wsx.png'  # L13 (stimulus)
$ws.Cells.Item(13, 13).Value = 79.07142857142857  # M13 (conceptual)
$ws.Cells.Item(13, 14).Value = 58  # N13 (perceptual)
$ws.Cells.Item(13, 15).Value = 68.53571428571428  # O13 (typicality)
$ws.Cells.Item(13, 16).Value = 42  # P13 (n)
$ws.Cells.Item(13, 17).Value = 7  # Q13 (p_typicality)
$ws.Cells.Item(13, 18).Value = 7  # R13 (p_conceptual)
$ws.Cells.Item(13, 19).Value = 7  # S13 (p_perceptual)
$ws.Cells.Item(13, 20).Value = 7  # T13 (r_typicality)
$ws.Cells.Item(13, 21).Value = 7  # U13 (r_conceptual)
$ws.Cells.Item(13, 22).Value = 7  # V13 (r_perceptual)

# Row 14
$ws.Cells.Item(14, 3).Value = 6  # C14 (block_total)
$ws.Cells.Item(14, 6).Value = 215  # F14 (trial_total)
$ws.Cells.Item(14, 9).Value = $null  # I14 (cond_cat)
$ws.Cells.Item(14, 10).Value = 'new'  # J14 (cond_mem)
$ws.Cells.Item(14, 11).Value = 'f'  # K14 (correct_answer)
$ws.Cells.Item(14, 12).Value = 'stimuli/img_vgaye.png'  # L14 (stimulus)
$ws.Cells.Item(14, 13).Value = 80.33333333333333  # M14 (conceptual)
$ws.Cells.Item(14, 14).Value = 64.57777777777778  # N14 (perceptual)
$ws.Cells.Item(14, 15).Value = 72.45555555555555  # O14 (typicality)
$ws.Cells.Item(14, 16).Value = 45  # P14 (n)
$ws.Cells.Item(14, 17).Value = 8  # Q14 (p_typicality)
$ws.Cells.Item(14, 18).Value = 8  # R14 (p_conceptual)
$ws.Cells.Item(14, 19).Value = 8  # S14 (p_perceptual)
$ws.Cells.Item(14, 20).Value = 8  # T14 (r_typicality)
$ws.Cells.Item(14, 21).Value = 7  # U14 (r_conceptual)
$ws.Cells.Item(14, 22).Value = 8  # V14 (r_perceptual)

# Row 15
$ws.Cells.Item(15, 3).Value = 6  # C15 (block_total)
$ws.Cells.Item(15, 6).Value = 216  # F15 (trial_total)
$ws.Cells.Item(15, 9).Value = $null  # I15 (cond_cat)
$ws.Cells.Item(15, 10).Value = 'new'  # J15 (cond_mem)
$ws.Cells.Item(15, 11).Value = 'f'  # K15 (correct_answer)
$ws.Cells.Item(15, 12).Value = 'stimuli/img_165pk.png'  # L15 (stimulus)
$ws.Cells.Item(15, 13).Value = 85.73333333333333  # M15 (conceptual)
$ws.Cells.Item(15, 14).Value = 69.22222222222223  # N15 (perceptual)
$ws.Cells.Item(15, 15).Value = 77.47777777777779  # O15 (typicality)
$ws.Cells.Item(15, 16).Value = 45  # P15 (n)
$ws.Cells.Item(15, 17).Value = 9  # Q15 (p_typicality)
$ws.Cells.Item(15, 18).Value = 9  # R15 (p_conceptual)
$ws.Cells.Item(15, 19).Value = 9  # S15 (p_perceptual)
$ws.Cells.Item(15, 20).Value = 9  # T15 (r_typicality)
$ws.Cells.Item(15, 21).Value = 9  # U15 (r_conceptual)
$ws.Cells.Item(15, 22).Value = 9  # V15 (r_perceptual)

# Row 16
$ws.Cells.Item(16, 3).Value = 6  # C16 (block_total)
$ws.Cells.Item(16, 6).Value = 217  # F16 (trial_total)
$ws.Cells.Item(16, 9).Value = $null  # I16 (cond_cat)
$ws.Cells.Item(16, 10).Value = 'new'  # J16 (cond_mem)
$ws.Cells.Item(16, 11).Value = 'f'  # K16 (correct_answer)
$ws.Cells.Item(16, 12).Value = 'stimuli/img_ra2nm.png'  # L16 (stimulus)
$ws.Cells.Item(16, 13).Value = 70.75  # M16 (conceptual)
$ws.Cells.Item(16, 14).Value = 50.375  # N16 (perceptual)
$ws.Cells.Item(16, 15).Value = 60.5625  # O16 (typicality)
$ws.Cells.Item(16, 16).Value = 40  # P16 (n)
$ws.Cells.Item(16, 17).Value = 6  # Q16 (p_typicality)
$ws.Cells.Item(16, 18).Value = 6  # R16 (p_conceptual)
$ws.Cells.Item(16, 19).Value = 6  # S16 (p_perceptual)
$ws.Cells.Item(16, 20).Value = 5  # T16 (r_typicality)
$ws.Cells.Item(16, 21).Value = 5  # U16 (r_conceptual)
$ws.Cells.Item(16, 22).Value = 6  # V16 (r_perceptual)

# Row 17
$ws.Cells.Item(17, 3).Value = 6  # C17 (block_total)
$ws.Cells.Item(17, 6).Value = 218  # F17 (trial_total)
$ws.Cells.Item(17, 9).Value = $null  # I17 (cond_cat)
$ws.Cells.Item(17, 10).Value = 'new'  # J17 (cond_mem)
$ws.Cells.Item(17, 11).Value = 'f'  # K17 (correct_answer)
$ws.Cells.Item(17, 12).Value = 'stimuli/img_xr3up.png'  # L17 (stimulus)
$ws.Cells.Item(17, 13).Value = 76.24444444444444  # M17 (conceptual)
$ws.Cells.Item(17, 14).Value = 55.88888888888889  # N17 (perceptual)
$ws.Cells.Item(17, 15).Value = 66.06666666666666  # O17 (typicality)
$ws.Cells.Item(17, 16).Value = 45  # P17 (n)
$ws.Cells.Item(17, 17).Value = 7  # Q17 (p_typicality)
$ws.Cells.Item(17, 18).Value = 7  # R17 (p_conceptual)
$ws.Cells.Item(17, 19).Value = 7  # S17 (p_perceptual)
$ws.Cells.Item(17, 20).Value = 6  # T17 (r_typicality)
$ws.Cells.Item(17, 21).Value = 6  # U17 (r_conceptual)
$ws.Cells.Item(17, 22).Value = 6  # V17 (r_perceptual)

# Row 18
$ws.Cells.Item(18, 3).Value = 6  # C18 (block_total)
$ws.Cells.Item(18, 6).Value = 219  # F18 (trial_total)
$ws.Cells.Item(18, 12).Value = 'stimuli/img_q9lab.png'  # L18 (stimulus)
$ws.Cells.Item(18, 13).Value = 53.97560975609756  # M18 (conceptual)
$ws.Cells.Item(18, 14).Value = 32.90243902439025  # N18 (perceptual)
$ws.Cells.Item(18, 15).Value = 43.4390243902439  # O18 (typicality)
$ws.Cells.Item(18, 16).Value = 41  # P18 (n)
$ws.Cells.Item(18, 17).Value = 3  # Q18 (p_typicality)
$ws.Cells.Item(18, 18).Value = 3  # R18 (p_conceptual)
$ws.Cells.Item(18, 19).Value = 3  # S18 (p_perceptual)
$ws.Cells.Item(18, 20).Value = 3  # T18 (r_typicality)
$ws.Cells.Item(18, 21).Value = 4  # U18 (r_conceptual)
$ws.Cells.Item(18, 22).Value = 3  # V18 (r_perceptual)

# Row 19
$ws.Cells.Item(19, 3).Value = 6  # C19 (block_total)
$ws.Cells.Item(19, 6).Value = 220  # F19 (trial_total)
$ws.Cells.Item(19, 12).Value = 'stimuli/img_pjfx6.png'  # L19 (stimulus)
$ws.Cells.Item(19, 13).Value = 32.23404255319149  # M19 (conceptual)
$ws.Cells.Item(19, 14).Value = 26.59574468085106  # N19 (perceptual)
$ws.Cells.Item(19, 15).Value = 29.41489361702127  # O19 (typicality)
$ws.Cells.Item(19, 16).Value = 47  # P19 (n)
$ws.Cells.Item(19, 17).Value = 2  # Q19 (p_typicality)
$ws.Cells.Item(19, 18).Value = 2  # R19 (p_conceptual)
$ws.Cells.Item(19, 19).Value = 2  # S19 (p_perceptual)
$ws.Cells.Item(19, 20).Value = 2  # T19 (r_typicality)
$ws.Cells.Item(19, 21).Value = 2  # U19 (r_conceptual)
$ws.Cells.Item(19, 22).Value = 3  # V19 (r_perceptual)

# Row 20
$ws.Cells.Item(20, 3).Value = 6  # C20 (block_total)
$ws.Cells.Item(20, 6).Value = 221  # F20 (trial_total)
$ws.Cells.Item(20, 12).Value = 'stimuli/img_gka64.png'  # L20 (stimulus)
$ws.Cells.Item(20, 13).Value = 19.23809523809524  # M20 (conceptual)
$ws.Cells.Item(20, 14).Value = 20.02380952380953  # N20 (perceptual)
$ws.Cells.Item(20, 15).Value = 19.63095238095238  # O20 (typicality)
$ws.Cells.Item(20, 16).Value = 42  # P20 (n)
$ws.Cells.Item(20, 17).Value = 1  # Q20 (p_typicality)
$ws.Cells.Item(20, 18).Value = 1  # R20 (p_conceptual)
$ws.Cells.Item(20, 19).Value = 1  # S20 (p_perceptual)
$ws.Cells.Item(20, 20).Value = 1  # T20 (r_typicality)
$ws.Cells.Item(20, 21).Value = 1  # U20 (r_conceptual)
$ws.Cells.Item(20, 22).Value = 2  # V20 (r_perceptual)

# Row 21
$ws.Cells.Item(21, 3).Value = 6  # C21 (block_total)
$ws.Cells.Item(21, 6).Value = 222  # F21 (trial_total)
$ws.Cells.Item(21, 12).Value = 'stimuli/img_c89x3.png'  # L21 (stimulus)
$ws.Cells.Item(21, 13).Value = 72.8695652173913  # M21 (conceptual)
$ws.Cells.Item(21, 14).Value = 49.65217391304348  # N21 (perceptual)
$ws.Cells.Item(21, 15).Value = 61.26086956521739  # O21 (typicality)
$ws.Cells.Item(21, 16).Value = 46  # P21 (n)
$ws.Cells.Item(21, 17).Value = 6  # Q21 (p_typicality)
$ws.Cells.Item(21, 18).Value = 6  # R21 (p_conceptual)
$ws.Cells.Item(21, 19).Value = 6  # S21 (p_perceptual)
$ws.Cells.Item(21, 21).Value = 6  # U21 (r_conceptual)
$ws.Cells.Item(21, 22).Value = 5  # V21 (r_perceptual)

# Row 22
$ws.Cells.Item(22, 3).Value = 6  # C22 (block_total)
$ws.Cells.Item(22, 6).Value = 223  # F22 (trial_total)
$ws.Cells.Item(22, 9).Value = $null  # I22 (cond_cat)
$ws.Cells.Item(22, 10).Value = 'new'  # J22 (cond_mem)
$ws.Cells.Item(22, 11).Value = 'f'  # K22 (correct_answer)
$ws.Cells.Item(22, 12).Value = 'stimuli/img_j856a.png'  # L22 (stimulus)
$ws.Cells.Item(22, 13).Value = 38.225  # M22 (conceptual)
$ws.Cells.Item(22, 14).Value = 25.875  # N22 (perceptual)
$ws.Cells.Item(22, 15).Value = 32.05  # O22 (typicality)
$ws.Cells.Item(22, 16).Value = 40  # P22 (n)
$ws.Cells.Item(22, 17).Value = 2  # Q22 (p_typicality)
$ws.Cells.Item(22, 18).Value = 2  # R22 (p_conceptual)
$ws.Cells.Item(22, 19).Value = 2  # S22 (p_perceptual)
$ws.Cells.Item(22, 20).Value = 3  # T22 (r_typicality)
$ws.Cells.Item(22, 21).Value = 3  # U22 (r_conceptual)
$ws.Cells.Item(22, 22).Value = 2  # V22 (r_perceptual)

# Row 23
$ws.Cells.Item(23, 3).Value = 6  # C23 (block_total)
$ws.Cells.Item(23, 6).Value = 224  # F23 (trial_total)
$ws.Cells.Item(23, 12).Value = 'stimuli/img_wbws6.png'  # L23 (stimulus)
$ws.Cells.Item(23, 13).Value = 57.97777777777777  # M23 (conceptual)
$ws.Cells.Item(23, 14).Value = 42.53333333333333  # N23 (perceptual)
$ws.Cells.Item(23, 15).Value = 50.25555555555555  # O23 (typicality)
$ws.Cells.Item(23, 16).Value = 45  # P23 (n)
$ws.Cells.Item(23, 17).Value = 4  # Q23 (p_typicality)
$ws.Cells.Item(23, 18).Value = 4  # R23 (p_conceptual)
$ws.Cells.Item(23, 19).Value = 4  # S23 (p_perceptual)
$ws.Cells.Item(23, 20).Value = 4  # T23 (r_typicality)
$ws.Cells.Item(23, 21).Value = 4  # U23 (r_conceptual)
$ws.Cells.Item(23, 22).Value = 5  # V23 (r_perceptual)

# Row 24
$ws.Cells.Item(24, 3).Value = 6  # C24 (block_total)
$ws.Cells.Item(24, 6).Value = 225  # F24 (trial_total)
$ws.Cells.Item(24, 9).Value = $null  # I24 (cond_cat)
$ws.Cells.Item(24, 10).Value = 'new'  # J24 (cond_mem)
$ws.Cells.Item(24, 11).Value = 'f'  # K24 (correct_answer)
$ws.Cells.Item(24, 12).Value = 'stimuli/img_g13d5.png'  # L24 (stimulus)
$ws.Cells.Item(24, 13).Value = 73  # M24 (conceptual)
$ws.Cells.Item(24, 14).Value = 51.51111111111111  # N24 (perceptual)
$ws.Cells.Item(24, 15).Value = 62.25555555555556  # O24 (typicality)
$ws.Cells.Item(24, 16).Value = 45  # P24 (n)
$ws.Cells.Item(24, 17).Value = 6  # Q24 (p_typicality)
$ws.Cells.Item(24, 18).Value = 6  # R24 (p_conceptual)
$ws.Cells.Item(24, 19).Value = 6  # S24 (p_perceptual)
$ws.Cells.Item(24, 20).Value = 6  # T24 (r_typicality)
$ws.Cells.Item(24, 21).Value = 6  # U24 (r_conceptual)
$ws.Cells.Item(24, 22).Value = 6  # V24 (r_perceptual)

# Row 25
$ws.Cells.Item(25, 3).Value = 6  # C25 (block_total)
$ws.Cells.Item(25, 6).Value = 226  # F25 (trial_total)
$ws.Cells.Item(25, 9).Value = 'target'  # I25 (cond_cat)
$ws.Cells.Item(25, 10).Value = 'old'  # J25 (cond_mem)
$ws.Cells.Item(25, 11).Value = 'j'  # K25 (correct_answer)
$ws.Cells.Item(25, 12).Value = 'stimuli/img_w8yhd.png'  # L25 (stimulus)
$ws.Cells.Item(25, 13).Value = 55.74418604651163  # M25 (conceptual)
$ws.Cells.Item(25, 14).Value = 38.90697674418605  # N25 (perceptual)
$ws.Cells.Item(25, 15).Value = 47.32558139534883  # O25 (typicality)
$ws.Cells.Item(25, 16).Value = 43  # P25 (n)
$ws.Cells.Item(25, 17).Value = 4  # Q25 (p_typicality)
$ws.Cells.Item(25, 18).Value = 4  # R25 (p_conceptual)
$ws.Cells.Item(25, 19).Value = 4  # S25 (p_perceptual)
$ws.Cells.Item(25, 20).Value = 4  # T25 (r_typicality)
$ws.Cells.Item(25, 22).Value = 4  # V25 (r_perceptual)

# Row 26
$ws.Cells.Item(26, 3).Value = 6  # C26 (block_total)
$ws.Cells.Item(26, 6).Value = 227  # F26 (trial_total)
$ws.Cells.Item(26, 8).Value = 'living_rooms'  # H26 (category)
$ws.Cells.Item(26, 10).Value = 'new'  # J26 (cond_mem)
$ws.Cells.Item(26, 12).Value = 'stimuli/img_b21d7.png'  # L26 (stimulus)
$ws.Cells.Item(26, 13).Value = 27.75555555555556  # M26 (conceptual)
$ws.Cells.Item(26, 14).Value = 13.86666666666667  # N26 (perceptual)
$ws.Cells.Item(26, 15).Value = 20.81111111111111  # O26 (typicality)
$ws.Cells.Item(26, 16).Value = 45  # P26 (n)
$ws.Cells.Item(26, 17).Value = 1  # Q26 (p_typicality)
$ws.Cells.Item(26, 18).Value = 1  # R26 (p_conceptual)
$ws.Cells.Item(26, 19).Value = 1  # S26 (p_perceptual)
$ws.Cells.Item(26, 20).Value = 2  # T26 (r_typicality)
$ws.Cells.Item(26, 21).Value = 2  # U26 (r_conceptual)
$ws.Cells.Item(26, 22).Value = 1  # V26 (r_perceptual)

# Row 27
$ws.Cells.Item(27, 3).Value = 6  # C27 (block_total)
$ws.Cells.Item(27, 6).Value = 228  # F27 (trial_total)
$ws.Cells.Item(27, 9).Value = 'target'  # I27 (cond_cat)
$ws.Cells.Item(27, 10).Value = 'old'  # J27 (cond_mem)
$ws.Cells.Item(27, 11).Value = 'j'  # K27 (correct_answer)
$ws.Cells.Item(27, 12).Value = 'stimuli/img_a9he3.png'  # L27 (stimulus)
$ws.Cells.Item(27, 13).Value = 83.06521739130434  # M27 (conceptual)
$ws.Cells.Item(27, 14).Value = 63.95652173913044  # N27 (perceptual)
$ws.Cells.Item(27, 15).Value = 73.51086956521739  # O27 (typicality)
$ws.Cells.Item(27, 16).Value = 46  # P27 (n)
$ws.Cells.Item(27, 17).Value = 8  # Q27 (p_typicality)
$ws.Cells.Item(27, 18).Value = 8  # R27 (p_conceptual)
$ws.Cells.Item(27, 19).Value = 8  # S27 (p_perceptual)
$ws.Cells.Item(27, 20).Value = 8  # T27 (r_typicality)
$ws.Cells.Item(27, 21).Value = 8  # U27 (r_conceptual)
$ws.Cells.Item(27, 22).Value = 8  # V27 (r_perceptual)

# Row 28
$ws.Cells.Item(28, 3).Value = 6  # C28 (block_total)
$ws.Cells.Item(28, 6).Value = 229  # F28 (trial_total)
$ws.Cells.Item(28, 8).Value = $null  # H28 (category)
$ws.Cells.Item(28, 9).Value = $null  # I28 (cond_cat)
$ws.Cells.Item(28, 10).Value = 'catch'  # J28 (cond_mem)
$ws.Cells.Item(28, 11).Value = 'f'  # K28 (correct_answer)
$ws.Cells.Item(28, 12).Value = 'stimuli/catch_25.jpg'  # L28 (stimulus)
$ws.Cells.Item(28, 13).Value = $null  # M28 (conceptual)
$ws.Cells.Item(28, 14).Value = $null  # N28 (perceptual)
$ws.Cells.Item(28, 15).Value = $null  # O28 (typicality)
$ws.Cells.Item(28, 16).Value = $null  # P28 (n)
$ws.Cells.Item(28, 17).Value = $null  # Q28 (p_typicality)
$ws.Cells.Item(28, 18).Value = $null  # R28 (p_conceptual)
$ws.Cells.Item(28, 19).Value = $null  # S28 (p_perceptual)
$ws.Cells.Item(28, 20).Value = $null  # T28 (r_typicality)
$ws.Cells.Item(28, 21).Value = $null  # U28 (r_conceptual)
$ws.Cells.Item(28, 22).Value = $null  # V28 (r_perceptual)

# Row 29
$ws.Cells.Item(29, 3).Value = 6  # C29 (block_total)
$ws.Cells.Item(29, 6).Value = 230  # F29 (trial_total)
$ws.Cells.Item(29, 9).Value = $null  # I29 (cond_cat)
$ws.Cells.Item(29, 10).Value = 'new'  # J29 (cond_mem)
$ws.Cells.Item(29, 11).Value = 'f'  # K29 (correct_answer)
$ws.Cells.Item(29, 12).Value = 'stimuli/img_tn8ys.png'  # L29 (stimulus)
$ws.Cells.Item(29, 13).Value = 86.70454545454545  # M29 (conceptual)
$ws.Cells.Item(29, 14).Value = 72.4090909090909  # N29 (perceptual)
$ws.Cells.Item(29, 15).Value = 79.55681818181819  # O29 (typicality)
$ws.Cells.Item(29, 16).Value = 44  # P29 (n)
$ws.Cells.Item(29, 17).Value = 10  # Q29 (p_typicality)
$ws.Cells.Item(29, 18).Value = 10  # R29 (p_conceptual)
$ws.Cells.Item(29, 19).Value = 10  # S29 (p_perceptual)
$ws.Cells.Item(29, 20).Value = 9  # T29 (r_typicality)
$ws.Cells.Item(29, 21).Value = 9  # U29 (r_conceptual)
$ws.Cells.Item(29, 22).Value = 10  # V29 (r_perceptual)

# Row 30
$ws.Cells.Item(30, 3).Value = 6  # C30 (block_total)
$ws.Cells.Item(30, 6).Value = 231  # F30 (trial_total)
$ws.Cells.Item(30, 9).Value = $null  # I30 (cond_cat)
$ws.Cells.Item(30, 10).Value = 'new'  # J30 (cond_mem)
$ws.Cells.Item(30, 11).Value = 'f'  # K30 (correct_answer)
$ws.Cells.Item(30, 12).Value = 'stimuli/img_jpldg.png'  # L30 (stimulus)
$ws.Cells.Item(30, 13).Value = 79.54545454545455  # M30 (conceptual)
$ws.Cells.Item(30, 14).Value = 57.75  # N30 (perceptual)
$ws.Cells.Item(30, 15).Value = 68.64772727272728  # O30 (typicality)
$ws.Cells.Item(30, 16).Value = 44  # P30 (n)
$ws.Cells.Item(30, 17).Value = 7  # Q30 (p_typicality)
$ws.Cells.Item(30, 18).Value = 7  # R30 (p_conceptual)
$ws.Cells.Item(30, 19).Value = 7  # S30 (p_perceptual)
$ws.Cells.Item(30, 20).Value = 7  # T30 (r_typicality)
$ws.Cells.Item(30, 21).Value = 7  # U30 (r_conceptual)
$ws.Cells.Item(30, 22).Value = 7  # V30 (r_perceptual)

# Row 31
$ws.Cells.Item(31, 3).Value = 6  # C31 (block_total)
$ws.Cells.Item(31, 6).Value = 232  # F31 (trial_total)
$ws.Cells.Item(31, 9).Value = 'target'  # I31 (cond_cat)
$ws.Cells.Item(31, 10).Value = 'old'  # J31 (cond_mem)
$ws.Cells.Item(31, 11).Value = 'j'  # K31 (correct_answer)
$ws.Cells.Item(31, 12).Value = 'stimuli/img_c0vzo.png'  # L31 (stimulus)
$ws.Cells.Item(31, 13).Value = 21.51162790697675  # M31 (conceptual)
$ws.Cells.Item(31, 14).Value = 8.232558139534884  # N31 (perceptual)
$ws.Cells.Item(31, 15).Value = 14.87209302325581  # O31 (typicality)
$ws.Cells.Item(31, 16).Value = 43  # P31 (n)
$ws.Cells.Item(31, 17).Value = 1  # Q31 (p_typicality)
$ws.Cells.Item(31, 18).Value = 1  # R31 (p_conceptual)
$ws.Cells.Item(31, 19).Value = 1  # S31 (p_perceptual)
$ws.Cells.Item(31, 20).Value = 1  # T31 (r_typicality)
$ws.Cells.Item(31, 21).Value = 1  # U31 (r_conceptual)
$ws.Cells.Item(31, 22).Value = 1  # V31 (r_perceptual)

# Row 32
$ws.Cells.Item(32, 3).Value = 6  # C32 (block_total)
$ws.Cells.Item(32, 6).Value = 233  # F32 (trial_total)
$ws.Cells.Item(32, 9).Value = 'target'  # I32 (cond_cat)
$ws.Cells.Item(32, 10).Value = 'old'  # J32 (cond_mem)
$ws.Cells.Item(32, 11).Value = 'j'  # K32 (correct_answer)
$ws.Cells.Item(32, 12).Value = 'stimuli/img_xzyzy.png'  # L32 (stimulus)
$ws.Cells.Item(32, 13).Value = 85.37209302325581  # M32 (conceptual)
$ws.Cells.Item(32, 14).Value = 68.90697674418605  # N32 (perceptual)
$ws.Cells.Item(32, 15).Value = 77.13953488372093  # O32 (typicality)

# Row 33
$ws.Cells.Item(33, 3).Value = 6  # C33 (block_total)
$ws.Cells.Item(33, 6).Value = 234  # F33 (trial_total)
$ws.Cells.Item(33, 9).Value = $null  # I33 (cond_cat)
$ws.Cells.Item(33, 10).Value = 'new'  # J33 (cond_mem)
$ws.Cells.Item(33, 11).Value = 'f'  # K33 (correct_answer)
$ws.Cells.Item(33, 12).Value = 'stimuli/img_3sw8t.png'  # L33 (stimulus)
$ws.Cells.Item(33, 13).Value = 67.4888888888889  # M33 (conceptual)
$ws.Cells.Item(33, 14).Value = 48.51111111111111  # N33 (perceptual)
$ws.Cells.Item(33, 15).Value = 58  # O33 (typicality)
$ws.Cells.Item(33, 16).Value = 45  # P33 (n)
$ws.Cells.Item(33, 17).Value = 5  # Q33 (p_typicality)
$ws.Cells.Item(33, 18).Value = 5  # R33 (p_conceptual)
$ws.Cells.Item(33, 19).Value = 5  # S33 (p_perceptual)
$ws.Cells.Item(33, 20).Value = 5  # T33 (r_typicality)
$ws.Cells.Item(33, 21).Value = 5  # U33 (r_conceptual)
$ws.Cells.Item(33, 22).Value = 5  # V33 (r_perceptual)

# Row 34
$ws.Cells.Item(34, 3).Value = 6  # C34 (block_total)
$ws.Cells.Item(34, 6).Value = 235  # F34 (trial_total)
$ws.Cells.Item(34, 9).Value = $null  # I34 (cond_cat)
$ws.Cells.Item(34, 10).Value = 'new'  # J34 (cond_mem)
$ws.Cells.Item(34, 11).Value = 'f'  # K34 (correct_answer)
$ws.Cells.Item(34, 12).Value = 'stimuli/img_x4bln.png'  # L34 (stimulus)
$ws.Cells.Item(34, 13).Value = 76.34042553191489  # M34 (conceptual)
$ws.Cells.Item(34, 14).Value = 59.51063829787234  # N34 (perceptual)
$ws.Cells.Item(34, 15).Value = 67.92553191489361  # O34 (typicality)
$ws.Cells.Item(34, 16).Value = 47  # P34 (n)
$ws.Cells.Item(34, 17).Value = 7  # Q34 (p_typicality)
$ws.Cells.Item(34, 18).Value = 7  # R34 (p_conceptual)
$ws.Cells.Item(34, 19).Value = 7  # S34 (p_perceptual)
$ws.Cells.Item(34, 20).Value = 7  # T34 (r_typicality)
$ws.Cells.Item(34, 21).Value = 7  # U34 (r_conceptual)
$ws.Cells.Item(34, 22).Value = 7  # V34 (r_perceptual)

# Row 35
$ws.Cells.Item(35, 3).Value = 6  # C35 (block_total)
$ws.Cells.Item(35, 6).Value = 236  # F35 (trial_total)
$ws.Cells.Item(35, 9).Value = 'target'  # I35 (cond_cat)
$ws.Cells.Item(35, 10).Value = 'old'  # J35 (cond_mem)
$ws.Cells.Item(35, 11).Value = 'j'  # K35 (correct_answer)
$ws.Cells.Item(35, 12).Value = 'stimuli/img_wgddx.png'  # L35 (stimulus)
$ws.Cells.Item(35, 13).Value = 45.6304347826087  # M35 (conceptual)
$ws.Cells.Item(35, 14).Value = 34.30434782608695  # N35 (perceptual)
$ws.Cells.Item(35, 15).Value = 39.96739130434783  # O35 (typicality)
$ws.Cells.Item(35, 16).Value = 46  # P35 (n)
$ws.Cells.Item(35, 20).Value = 3  # T35 (r_typicality)
$ws.Cells.Item(35, 22).Value = 3  # V35 (r_perceptual)

# Row 36
$ws.Cells.Item(36, 3).Value = 6  # C36 (block_total)
$ws.Cells.Item(36, 6).Value = 237  # F36 (trial_total)
$ws.Cells.Item(36, 9).Value = 'target'  # I36 (cond_cat)
$ws.Cells.Item(36, 10).Value = 'old'  # J36 (cond_mem)
$ws.Cells.Item(36, 11).Value = 'j'  # K36 (correct_answer)
$ws.Cells.Item(36, 12).Value = 'stimuli/img_pbsj1.png'  # L36 (stimulus)
$ws.Cells.Item(36, 13).Value = 73.88636363636364  # M36 (conceptual)
$ws.Cells.Item(36, 14).Value = 51.52272727272727  # N36 (perceptual)
$ws.Cells.Item(36, 15).Value = 62.70454545454545  # O36 (typicality)
$ws.Cells.Item(36, 16).Value = 44  # P36 (n)
$ws.Cells.Item(36, 17).Value = 6  # Q36 (p_typicality)
$ws.Cells.Item(36, 18).Value = 6  # R36 (p_conceptual)
$ws.Cells.Item(36, 19).Value = 6  # S36 (p_perceptual)
$ws.Cells.Item(36, 20).Value = 6  # T36 (r_typicality)
$ws.Cells.Item(36, 21).Value = 6  # U36 (r_conceptual)
$ws.Cells.Item(36, 22).Value = 6  # V36 (r_perceptual)

# Row 37
$ws.Cells.Item(37, 3).Value = 6  # C37 (block_total)
$ws.Cells.Item(37, 6).Value = 238  # F37 (trial_total)
$ws.Cells.Item(37, 9).Value = 'target'  # I37 (cond_cat)
$ws.Cells.Item(37, 10).Value = 'old'  # J37 (cond_mem)
$ws.Cells.Item(37, 11).Value = 'j'  # K37 (correct_answer)
$ws.Cells.Item(37, 12).Value = 'stimuli/img_8dmpq.png'  # L37 (stimulus)
$ws.Cells.Item(37, 13).Value = 30.65909090909091  # M37 (conceptual)
$ws.Cells.Item(37, 14).Value = 24.11363636363636  # N37 (perceptual)
$ws.Cells.Item(37, 15).Value = 27.38636363636364  # O37 (typicality)
$ws.Cells.Item(37, 16).Value = 44  # P37 (n)
$ws.Cells.Item(37, 17).Value = 2  # Q37 (p_typicality)
$ws.Cells.Item(37, 18).Value = 2  # R37 (p_conceptual)
$ws.Cells.Item(37, 19).Value = 2  # S37 (p_perceptual)
$ws.Cells.Item(37, 20).Value = 2  # T37 (r_typicality)
$ws.Cells.Item(37, 21).Value = 2  # U37 (r_conceptual)
$ws.Cells.Item(37, 22).Value = 2  # V37 (r_perceptual)

# Row 38
$ws.Cells.Item(38, 3).Value = 6  # C38 (block_total)
$ws.Cells.Item(38, 6).Value = 239  # F38 (trial_total)
$ws.Cells.Item(38, 9).Value = 'target'  # I38 (cond_cat)
$ws.Cells.Item(38, 10).Value = 'old'  # J38 (cond_mem)
$ws.Cells.Item(38, 11).Value = 'j'  # K38 (correct_answer)
$ws.Cells.Item(38, 12).Value = 'stimuli/img_rych7.png'  # L38 (stimulus)
$ws.Cells.Item(38, 13).Value = 30.4468085106383  # M38 (conceptual)
$ws.Cells.Item(38, 14).Value = 23.4468085106383  # N38 (perceptual)
$ws.Cells.Item(38, 15).Value = 26.9468085106383  # O38 (typicality)
$ws.Cells.Item(38, 16).Value = 47  # P38 (n)
$ws.Cells.Item(38, 17).Value = 2  # Q38 (p_typicality)
$ws.Cells.Item(38, 18).Value = 2  # R38 (p_conceptual)
$ws.Cells.Item(38, 19).Value = 2  # S38 (p_perceptual)
$ws.Cells.Item(38, 20).Value = 2  # T38 (r_typicality)
$ws.Cells.Item(38, 21).Value = 2  # U38 (r_conceptual)
$ws.Cells.Item(38, 22).Value = 2  # V38 (r_perceptual)

# Row 39
$ws.Cells.Item(39, 3).Value = 6  # C39 (block_total)
$ws.Cells.Item(39, 6).Value = 240  # F39 (trial_total)
$ws.Cells.Item(39, 9).Value = 'target'  # I39 (cond_cat)
$ws.Cells.Item(39, 10).Value = 'old'  # J39 (cond_mem)
$ws.Cells.Item(39, 11).Value = 'j'  # K39 (correct_answer)
$ws.Cells.Item(39, 12).Value = 'stimuli/img_dg5h7.png'  # L39 (stimulus)
$ws.Cells.Item(39, 13).Value = 88.72093023255815  # M39 (conceptual)
$ws.Cells.Item(39, 14).Value = 76.06976744186046  # N39 (perceptual)
$ws.Cells.Item(39, 15).Value = 82.3953488372093  # O39 (typicality)
$ws.Cells.Item(39, 17).Value = 10  # Q39 (p_typicality)
$ws.Cells.Item(39, 18).Value = 10  # R39 (p_conceptual)
$ws.Cells.Item(39, 19).Value = 10  # S39 (p_perceptual)
$ws.Cells.Item(39, 20).Value = 10  # T39 (r_typicality)
$ws.Cells.Item(39, 21).Value = 10  # U39 (r_conceptual)
$ws.Cells.Item(39, 22).Value = 10  # V39 (r_perceptual)

# Row 40
$ws.Cells.Item(40, 3).Value = 6  # C40 (block_total)
$ws.Cells.Item(40, 6).Value = 241  # F40 (trial_total)
$ws.Cells.Item(40, 9).Value = 'target'  # I40 (cond_cat)
$ws.Cells.Item(40, 10).Value = 'old'  # J40 (cond_mem)
$ws.Cells.Item(40, 11).Value = 'j'  # K40 (correct_answer)
$ws.Cells.Item(40, 12).Value = 'stimuli/img_nb8p4.png'  # L40 (stimulus)
$ws.Cells.Item(40, 13).Value = 16.36170212765957  # M40 (conceptual)
$ws.Cells.Item(40, 14).Value = 12.70212765957447  # N40 (perceptual)
$ws.Cells.Item(40, 15).Value = 14.53191489361702  # O40 (typicality)
$ws.Cells.Item(40, 16).Value = 47  # P40 (n)
$ws.Cells.Item(40, 20).Value = 1  # T40 (r_typicality)
$ws.Cells.Item(40, 22).Value = 1  # V40 (r_perceptual)

# Row 41
$ws.Cells.Item(41, 3).Value = 6  # C41 (block_total)
$ws.Cells.Item(41, 6).Value = 242  # F41 (trial_total)
$ws.Cells.Item(41, 12).Value = 'stimuli/img_hmmra.png'  # L41 (stimulus)
$ws.Cells.Item(41, 13).Value = 54.65853658536585  # M41 (conceptual)
$ws.Cells.Item(41, 14).Value = 34.24390243902439  # N41 (perceptual)
$ws.Cells.Item(41, 15).Value = 44.45121951219512  # O41 (typicality)
$ws.Cells.Item(41, 16).Value = 41  # P41 (n)
$ws.Cells.Item(41, 17).Value = 3  # Q41 (p_typicality)
$ws.Cells.Item(41, 18).Value = 3  # R41 (p_conceptual)
$ws.Cells.Item(41, 19).Value = 3  # S41 (p_perceptual)
$ws.Cells.Item(41, 20).Value = 4  # T41 (r_typicality)
$ws.Cells.Item(41, 21).Value = 4  # U41 (r_conceptual)
$ws.Cells.Item(41, 22).Value = 3  # V41 (r_perceptual)

# Row 42
$ws.Cells.Item(42, 3).Value = 6  # C42 (block_total)
$ws.Cells.Item(42, 6).Value = 243  # F42 (trial_total)
$ws.Cells.Item(42, 9).Value = 'target'  # I42 (cond_cat)
$ws.Cells.Item(42, 10).Value = 'old'  # J42 (cond_mem)
$ws.Cells.Item(42, 11).Value = 'j'  # K42 (correct_answer)
$ws.Cells.Item(42, 12).Value = 'stimuli/img_zxvl3.png'  # L42 (stimulus)
$ws.Cells.Item(42, 13).Value = 68.78260869565217  # M42 (conceptual)
$ws.Cells.Item(42, 14).Value = 47.56521739130435  # N42 (perceptual)
$ws.Cells.Item(42, 15).Value = 58.17391304347827  # O42 (typicality)
$ws.Cells.Item(42, 16).Value = 46  # P42 (n)
$ws.Cells.Item(42, 17).Value = 5  # Q42 (p_typicality)
$ws.Cells.Item(42, 18).Value = 5  # R42 (p_conceptual)
$ws.Cells.Item(42, 19).Value = 5  # S42 (p_perceptual)
$ws.Cells.Item(42, 20).Value = 5  # T42 (r_typicality)
$ws.Cells.Item(42, 21).Value = 5  # U42 (r_conceptual)
$ws.Cells.Item(42, 22).Value = 5  # V42 (r_perceptual)

